$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.626.05"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "'1.695.69"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'315.00"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.3919"
$ws.Range("E7").Value = "  -0.91%  "
$ws.Range("D8").Value = "'0.4024"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'1.517"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "'1.001"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "'53.00"
$ws.Range("E11").Value = "  -2.35%  "
$ws.Range("D12").Value = "'0.08829"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "'7.468"
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("D14").Value = "'23.57"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "'8.171"
$ws.Range("E15").Value = "  +7.65%  "
$ws.Range("D16").Value = "'0.00001317"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "'1.701.35"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").Value = "'99.42"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").Value = "'0.07017"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'19.60"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "'7.061"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("D22").Value = "'1.005"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "'14.63"
$ws.Range("E23").Value = "  +4.01%  "
$ws.Range("D24").Value = "'24.631.96"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "'3.118"
$ws.Range("E25").Value = "  +2.57%  "
$ws.Range("D26").Value = "'2.362"
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("D27").Value = "'22.57"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("D28").Value = "'162.59"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").Value = "'8.712"
$ws.Range("E29").Value = "  +15.20%  "
$ws.Range("D30").Value = "'135.38"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").Value = "'5.142"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").Value = "'0.08953"
$ws.Range("E32").Value = "  +4.75%  "
$ws.Range("D33").Value = "'7.582"
$ws.Range("E33").Value = "  +3.23%  "
$ws.Range("D34").Value = "'1.062"
$ws.Range("E34").Value = "  -3.38%  "
$ws.Range("D35").Value = "'1.962"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("D36").Value = "'11.03"
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("D37").Value = "'0.2748"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").Value = "'14.40"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").Value = "'0.02781"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("D40").Value = "'0.09095"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").Value = "'1.455"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").Value = "'0.7653"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").Value = "'15.91"
$ws.Range("E43").Value = "  +3.49%  "
$ws.Range("D44").Value = "'0.7142"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "'2.542"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").Value = "'4.209"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").Value = "'1.340"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").Value = "'139.54"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("D50").Value = "'0.07962"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").Value = "'89.97"
$ws.Range("E51").Value = "  +1.79%  "
